# Auto-generated Excel COM-interop script to update cryptos list values
# per commit: "Updated cryptos list on Wed Nov 15 16:33:21 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.362.48"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.014.09"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.49"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.20"
$ws.Range("E7").Value = "  +10.03%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.85"
$ws.Range("E9").Value = "  -8.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.370"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.901"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.81"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "2.309.86"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.42"
$ws.Range("E16").Value = "  +14.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.44"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "1.996.35"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "36.329.34"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.98"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.29"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.11"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.67"
$ws.Range("E24").Value = "  +18.24%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.07"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.60"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("E32").Value = "  +26.54%  "
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +10.33%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.86"
$ws.Range("E39").Value = "  +16.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +15.25%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.75"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.446.76"
$ws.Range("E46").Value = "  +4.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.77"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.82"
$ws.Range("E48").Value = "  +4.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").Value = "  +13.96%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.10"
$ws.Range("E51").Value = "  +2.72%  "
